# Generate Report for Handoff
#
# The "29849ffe..." file finished its handoff since the report was last
# generated: its priority is no longer "low" (it's now "ht" - high/hot) and
# its handoff timestamp moved forward a few seconds. Reflect that refreshed
# report state across all three sheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# zh-cn sheet: rows 4-7 correspond to 29849ffe, 2bb4b73c, f125352b, f20c9490
$zhcn = $wb.Worksheets.Item("zh-cn")
for ($r = 4; $r -le 7; $r++) {
  $zhcn.Cells.Item($r, 5).Value = "ht"                        # E: Priority
  $zhcn.Cells.Item($r, 8).Value = "2016-08-28 22:32:24"        # H: Latest Handoff Datetime
}

# de-de sheet: same rows get the refreshed priority, and the Latest Handoff
# Datetime column (which shares its text with the Overview generation date
# below) moves to the new generation timestamp as well
$dede = $wb.Worksheets.Item("de-de")
for ($r = 4; $r -le 7; $r++) {
  $dede.Cells.Item($r, 5).Value = "ht"                        # E: Priority
  $dede.Cells.Item($r, 8).Value = "2016-08-28 22:32:28"       # H: Latest Handoff Datetime
}

# Overview sheet: the "Latest HO Xliff Generate Date" column for the same
# files advances to the new generation timestamp
$overview = $wb.Worksheets.Item("Overview")
for ($r = 4; $r -le 7; $r++) {
  $overview.Cells.Item($r, 7).Value = "2016-08-28 22:32:28"   # G: Latest HO Xliff Generate Date
}
